# Append the 2025 poker year-end standings to the bottom of the data table.
# (Sheet1 held yearly rows through 2024 in rows 2-200; this adds the new
# 2025 season rows 201-210, one per player, matching the existing layout:
# Yr | Person | SRank | Points | Bonus | PointsBonus | Chips | Winnings |
# Takehome | PersStatus | pers_personid)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 201; Person = "Andy";     SRank = 1;  Points = 12; Bonus = 0; PointsBonus = 12; Chips = 44600; Winnings = 60; Takehome = 40;  PersonId = 349 },
    @{ Row = 202; Person = "Prashant"; SRank = 2;  Points = 9;  Bonus = 0; PointsBonus = 9;  Chips = 27800; Winnings = 20; Takehome = 0;   PersonId = 365 },
    @{ Row = 203; Person = "Matt";     SRank = 3;  Points = 9;  Bonus = 0; PointsBonus = 9;  Chips = 27050; Winnings = 50; Takehome = 30;  PersonId = 362 },
    @{ Row = 204; Person = "Richard";  SRank = 4;  Points = 6;  Bonus = 0; PointsBonus = 6;  Chips = 20800; Winnings = 20; Takehome = 0;   PersonId = 366 },
    @{ Row = 205; Person = "Pepe";     SRank = 5;  Points = 6;  Bonus = 0; PointsBonus = 6;  Chips = 20250; Winnings = 0;  Takehome = -20; PersonId = 364 },
    @{ Row = 206; Person = "Maisy";    SRank = 6;  Points = 5;  Bonus = 0; PointsBonus = 5;  Chips = 18350; Winnings = 0;  Takehome = -20; PersonId = 360 },
    @{ Row = 207; Person = "Mark";     SRank = 7;  Points = 5;  Bonus = 0; PointsBonus = 5;  Chips = 17050; Winnings = 10; Takehome = 0;   PersonId = 361 },
    @{ Row = 208; Person = "Jon";      SRank = 8;  Points = 4;  Bonus = 0; PointsBonus = 4;  Chips = 13050; Winnings = 0;  Takehome = -20; PersonId = 357 },
    @{ Row = 209; Person = "Anthony";  SRank = 9;  Points = 4;  Bonus = 0; PointsBonus = 4;  Chips = 9300;  Winnings = 10; Takehome = 0;   PersonId = 350 },
    @{ Row = 210; Person = "Alex";     SRank = 10; Points = 0;  Bonus = 0; PointsBonus = 0;  Chips = 4500;  Winnings = 0;  Takehome = -10; PersonId = 348 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Range("A$r").Value = 2025
    $ws.Range("B$r").Value = $entry.Person
    $ws.Range("C$r").Value = $entry.SRank
    $ws.Range("D$r").Value = $entry.Points
    $ws.Range("E$r").Value = $entry.Bonus
    $ws.Range("F$r").Value = $entry.PointsBonus
    $ws.Range("G$r").Value = $entry.Chips
    $ws.Range("H$r").Value = $entry.Winnings
    $ws.Range("I$r").Value = $entry.Takehome
    $ws.Range("J$r").Value = "Active"
    $ws.Range("K$r").Value = $entry.PersonId
}

# Reselect the full data range (now through row 210) like the source
# workbook's saved selection.
[void]$ws.Range("A1:K210").Select()
